$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.2457683333333333
$ws.Cells.Item(2, 8).Value = 0.737305
$ws.Cells.Item(2, 9).Value = 0.1447271191911903
$ws.Cells.Item(2, 10).Value = 0.1575855905380038
$ws.Cells.Item(2, 13).Value = 1.684496
$ws.Cells.Item(2, 14).Value = 5.053488
$ws.Cells.Item(2, 15).Value = 0.6423607101334534
$ws.Cells.Item(2, 16).Value = 0.7282461611889918
$ws.Cells.Item(2, 17).Value = 0.4139957744266667
$ws.Cells.Item(2, 18).Value = 3.72596196984
$ws.Cells.Item(2, 19).Value = 0.09296701505922197
$ws.Cells.Item(2, 20).Value = 0.1147611013680015
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.2457683333333333
$ws.Cells.Item(3, 8).Value = 0.737305
$ws.Cells.Item(3, 9).Value = 0.1447271191911903
$ws.Cells.Item(3, 10).Value = 0.1575855905380038
$ws.Cells.Item(3, 13).Value = 0.9277985
$ws.Cells.Item(3, 14).Value = 1.855597
$ws.Cells.Item(3, 15).Value = 0.3538039290807178
$ws.Cells.Item(3, 16).Value = 0.2674056793968462
$ws.Cells.Item(3, 17).Value = 0.2280234910141667
$ws.Cells.Item(3, 18).Value = 1.368140946085
$ws.Cells.Item(3, 19).Value = 0.0512050234143765
$ws.Cells.Item(3, 20).Value = 0.04213928190096811
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.2457683333333333
$ws.Cells.Item(4, 8).Value = 0.737305
$ws.Cells.Item(4, 9).Value = 0.1447271191911903
$ws.Cells.Item(4, 10).Value = 0.1575855905380038
$ws.Cells.Item(4, 13).Value = 0.01005766666666667
$ws.Cells.Item(4, 14).Value = 0.030173
$ws.Cells.Item(4, 15).Value = 0.003835360785828855
$ws.Cells.Item(4, 16).Value = 0.004348159414162149
$ws.Cells.Item(4, 17).Value = 0.002471855973888889
$ws.Cells.Item(4, 18).Value = 0.022246703765
$ws.Cells.Item(4, 19).Value = 0.0005550807175918702
$ws.Cells.Item(4, 20).Value = 0.0006852072690341226
$ws.Cells.Item(5, 9).Value = 0.6104826439049008
$ws.Cells.Item(5, 10).Value = 0.6647217777192627
$ws.Cells.Item(5, 13).Value = 1.684496
$ws.Cells.Item(5, 14).Value = 5.053488
$ws.Cells.Item(5, 15).Value = 0.6423607101334534
$ws.Cells.Item(5, 16).Value = 0.7282461611889918
$ws.Cells.Item(5, 17).Value = 1.746301842736
$ws.Cells.Item(5, 18).Value = 15.716716584624
$ws.Cells.Item(5, 19).Value = 0.3921500646629003
$ws.Cells.Item(5, 20).Value = 0.4840810828827753
$ws.Cells.Item(6, 9).Value = 0.6104826439049008
$ws.Cells.Item(6, 10).Value = 0.6647217777192627
$ws.Cells.Item(6, 13).Value = 0.9277985
$ws.Cells.Item(6, 14).Value = 1.855597
$ws.Cells.Item(6, 15).Value = 0.3538039290807178
$ws.Cells.Item(6, 16).Value = 0.2674056793968462
$ws.Cells.Item(6, 17).Value = 0.9618403547635
$ws.Cells.Item(6, 18).Value = 5.771042128581
$ws.Cells.Item(6, 19).Value = 0.2159911580491387
$ws.Cells.Item(6, 20).Value = 0.1777503785808988
$ws.Cells.Item(7, 9).Value = 0.6104826439049008
$ws.Cells.Item(7, 10).Value = 0.6647217777192627
$ws.Cells.Item(7, 13).Value = 0.01005766666666667
$ws.Cells.Item(7, 14).Value = 0.030173
$ws.Cells.Item(7, 15).Value = 0.003835360785828855
$ws.Cells.Item(7, 16).Value = 0.004348159414162149
$ws.Cells.Item(7, 17).Value = 0.01042669251433333
$ws.Cells.Item(7, 18).Value = 0.09384023262899999
$ws.Cells.Item(7, 19).Value = 0.002341421192861977
$ws.Cells.Item(7, 20).Value = 0.002890316255588611
$ws.Cells.Item(8, 7).Value = 0.4156905
$ws.Cells.Item(8, 8).Value = 0.831381
$ws.Cells.Item(8, 9).Value = 0.2447902369039089
$ws.Cells.Item(8, 10).Value = 0.1776926317427335
$ws.Cells.Item(8, 13).Value = 1.684496
$ws.Cells.Item(8, 14).Value = 5.053488
$ws.Cells.Item(8, 15).Value = 0.6423607101334534
$ws.Cells.Item(8, 16).Value = 0.7282461611889918
$ws.Cells.Item(8, 17).Value = 0.7002289844880001
$ws.Cells.Item(8, 18).Value = 4.201373906928
$ws.Cells.Item(8, 19).Value = 0.1572436304113312
$ws.Cells.Item(8, 20).Value = 0.1294039769382148
$ws.Cells.Item(9, 7).Value = 0.4156905
$ws.Cells.Item(9, 8).Value = 0.831381
$ws.Cells.Item(9, 9).Value = 0.2447902369039089
$ws.Cells.Item(9, 10).Value = 0.1776926317427335
$ws.Cells.Item(9, 13).Value = 0.9277985
$ws.Cells.Item(9, 14).Value = 1.855597
$ws.Cells.Item(9, 15).Value = 0.3538039290807178
$ws.Cells.Item(9, 16).Value = 0.2674056793968462
$ws.Cells.Item(9, 17).Value = 0.38567702236425
$ws.Cells.Item(9, 18).Value = 1.542708089457
$ws.Cells.Item(9, 19).Value = 0.08660774761720269
$ws.Cells.Item(9, 20).Value = 0.04751601891497925
$ws.Cells.Item(10, 7).Value = 0.4156905
$ws.Cells.Item(10, 8).Value = 0.831381
$ws.Cells.Item(10, 9).Value = 0.2447902369039089
$ws.Cells.Item(10, 10).Value = 0.1776926317427335
$ws.Cells.Item(10, 13).Value = 0.01005766666666667
$ws.Cells.Item(10, 14).Value = 0.030173
$ws.Cells.Item(10, 15).Value = 0.003835360785828855
$ws.Cells.Item(10, 16).Value = 0.004348159414162149
$ws.Cells.Item(10, 17).Value = 0.0041808764855
$ws.Cells.Item(10, 18).Value = 0.025085258913
$ws.Cells.Item(10, 19).Value = 0.0009388588753750075
$ws.Cells.Item(10, 20).Value = 0.0007726358895394145